# Applies the "Education" section edits:
#   1) "2022     Bachelor of Computing Edinburgh Napier University"
#        -> "2022     Bachelor of Computing, Edinburgh Napier University"
#      (splits the run so the added ", " becomes its own run(s), leaving the
#       existing "Edinburgh Napier University" run(s) untouched)
#   2) "College diploma in Computer Science,"
#        -> "College Diploma in Computer Science,"
#      (splits the run so "College ", "D" and "iploma in Computer Science,"
#       become separate runs)
#
# Plain Range.Text / Find-replace assignment on this runtime rebuilds
# (coalesces) every run in the touched paragraph, which would lose the
# pre-existing run boundaries. Splitting the paragraph with
# InsertParagraphAfter() at the desired character position -- then deleting
# the newly made paragraph mark to rejoin -- inserts/divides text without
# disturbing sibling runs, so it is used here instead.

$d = $word.ActiveDocument

function Find-RangeFromStart($doc, $text) {
    $rng = $doc.Content
    $ok = $rng.Find.Execute($text, $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for text: $text"
    }
    return $rng
}

function Split-AndInsert($doc, $pos, $insertText) {
    # Break the paragraph at $pos, type $insertText at the start of the new
    # paragraph (this becomes its own run since it starts a fresh paragraph),
    # then delete the paragraph mark to rejoin -- the inserted text ends up
    # as (a) separate run(s) ahead of whatever followed $pos, with no
    # formatting residue and no merging of later sibling runs.
    $doc.Range($pos, $pos).InsertParagraphAfter()
    if ($insertText.Length -gt 0) {
        $doc.Range($pos + 1, $pos + 1).InsertBefore($insertText)
    }
    $mark = $doc.Range($pos, $pos + 1)
    $mark.Delete()
}

# ---------------------------------------------------------------------
# Edit 1: "2022     Bachelor of Computing " -> "...Computing, " (comma)
# ---------------------------------------------------------------------
$r1 = Find-RangeFromStart $d "2022     Bachelor of Computing"
$pos1 = $r1.End
Split-AndInsert $d $pos1 ","

# ---------------------------------------------------------------------
# Edit 2: "College diploma in Computer Science," -> "College Diploma..."
#   a) split "College " from "d" from "iploma in Computer Science,"
#      (both breaks made before any rejoin, so the offsets of the later
#      break are not disturbed by the earlier one)
#   b) replace the isolated "d" run's text with "D"
#   c) rejoin the three temporary paragraphs back into one
# ---------------------------------------------------------------------
$r2 = Find-RangeFromStart $d "College diploma in Computer Science,"
$collegeEnd = $r2.Start + 8   # end of "College "
$dEnd = $collegeEnd + 1       # end of "d"

$d.Range($collegeEnd, $collegeEnd).InsertParagraphAfter()
$d.Range($dEnd + 1, $dEnd + 1).InsertParagraphAfter()

$dRange = $d.Range($collegeEnd + 1, $dEnd + 1)
$dRange.Text = "D"

# Rejoin: delete the mark ending "...College " (at $collegeEnd), then the
# mark ending the isolated "D" (originally at $dEnd+1, shifted down by 1
# once the first mark above is removed).
$d.Range($collegeEnd, $collegeEnd + 1).Delete()
$d.Range($dEnd, $dEnd + 1).Delete()
